$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("B10").Value = "Administrador acessa a funcionalidade de 'Catalogo (Perfis) de Competencias' a partir do menu inicial"
$ws.Range("D10").Value = "SYSTEM exibe a listagem do Catalogo (Perfis) de Competencias cadastradas com a opcao 'Alterar Gerente' dentre as varias listadas"
